$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, copying the formatting (font, border, alignment)
# from the neighboring header cell G1 so the new column matches the others.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H10 with 0 (plain numeric cells, no special style)
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
